# US 3.3 commit files
# Update the CIRbTF sheet: add three new fuel rows (heavy/residual fuel oil,
# LPG propane/butane, hydrogen) each flagged with -1, reword the "ratio or
# flag" header to "ratio or flag (dimensionless)" and wrap it, and widen
# column A to fit the longer fuel names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CIRbTF")

# New data rows first, so their strings are appended to the shared-string
# table ahead of the reworded header text.
$ws.Range("A9").Value = "heavy or residual fuel oil"
$ws.Range("B9").Value = -1

$ws.Range("A10").Value = "LPG propane or butane"
$ws.Range("B10").Value = -1

$ws.Range("A11").Value = "hydrogen"
$ws.Range("B11").Value = -1

# Header cell B1: new wording, wrapped text, taller row to fit the wrap
$ws.Range("B1").Value = "ratio or flag (dimensionless)"
$ws.Range("B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

# Column A a bit wider to fit the longer fuel names
$ws.Columns.Item(1).ColumnWidth = 24.8

# Leave CIRbTF's cursor on B1 (matching the saved view state), then hop back
# to the About tab so it stays the workbook's active sheet, same as before.
$ws.Range("B1").Select()
$wb.Worksheets.Item("About").Activate()
